$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtractedData")

$data = @(
    @("Opel Corsa D 1.2", "https://www.olx.ro/d/oferta/opel-corsa-d-1-2-IDfyj8M.html#714b7e881f;promoted", "2 700 €", "Timisoara", "Azi 22:52"),
    @("Opel Corsa AUTOMATA Benzina 1.2 Climatronic RATE", "https://www.autovit.ro/anunt/opel-corsa-1-2-ID7GWly1.html#xtor=SEC-81", "3 990 €", "Bacau", "Azi 22:16"),
    @("Opel Corsa 1,2 benzina Automata", "https://www.olx.ro/d/oferta/opel-corsa-1-2-benzina-automata-IDfvuA0.html#714b7e881f;promoted", "2 790 €", "Ghirdoveni", "Azi 21:58"),
    @("Opel Corsa 2004 / 1.0 benzina / 60cp / Stare generala buna/recent adus", "https://www.olx.ro/d/oferta/opel-corsa-2004-1-0-benzina-60cp-stare-generala-buna-recent-adus-IDfDuhE.html#714b7e881f;promoted", "1 299 €", "Ganeasa", "Azi 21:27"),
    @("Opel Corsa 2009/1.4 Benzina/Posibilitate Plata RATE", "https://www.olx.ro/d/oferta/opel-corsa-2009-1-4-benzina-posibilitate-plata-rate-IDfEZGA.html#714b7e881f;promoted", "3 299 €", "Brad", "Azi 21:14"),
    @("Opel Corsa 2007 Automata Benzina RATE", "https://www.autovit.ro/anunt/opel-corsa-1-2i-ID7GL1Mr.html#xtor=SEC-81", "3 990 €", "Bacau", "Azi 21:03"),
    @("Opel corsa 1.2 benzina", "https://www.olx.ro/d/oferta/opel-corsa-1-2-benzina-IDfrI9M.html#714b7e881f;promoted", "2 800 €", "Buzau", "Azi 20:45"),
    @("Vand Opel Corsa C", "https://www.olx.ro/d/oferta/vand-opel-corsa-c-IDfFFaV.html#714b7e881f", 707, "Vaslui", "Azi 20:35"),
    @("Opel Corsa 2011 AUTOMATA Euro 5 Rate", "https://www.autovit.ro/anunt/opel-corsa-ID7GJECP.html#xtor=SEC-81", "5 850 €", "Bacau", "Azi 20:24"),
    @("Opel corsa D 1.3 90cp", "https://www.olx.ro/d/oferta/opel-corsa-d-1-3-90cp-IDfFF42.html#714b7e881f", "1 750 €", "Brasov", "Azi 20:24")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$wb.Save()
